# [95] docs - add histogram chart
# Refresh the "Cartografia-social-api-front" repository rows at the tail of
# the analysis sheet: row 43 gets new (v1.0.2) figures, rows 44-51 are new
# version snapshots (v1.0.3 .. v2.0.1) appended after it, growing the sheet
# from A1:N45 to A1:N51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(43, 1).Value = 0.8043478260869565
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(43, 3).Value = 0.9565217391304348
$ws.Cells.Item(43, 4).Value = 0.9565217391304348
$ws.Cells.Item(43, 5).Value = 0.9565217391304348
$ws.Cells.Item(43, 6).Value = 0.9565217391304348
$ws.Cells.Item(43, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(43, 8).Value = "v1.0.2"
$ws.Cells.Item(43, 9).Value = 2657
$ws.Cells.Item(43, 10).Value = 0.5810869565217391
$ws.Cells.Item(43, 11).Value = 0.9565217391304348
$ws.Cells.Item(43, 12).Value = 0.2905434782608696
$ws.Cells.Item(43, 13).Value = 0.4782608695652174
$ws.Cells.Item(43, 14).Value = 0.768804347826087

$ws.Cells.Item(44, 1).Value = 0.8
$ws.Cells.Item(44, 2).Value = 0
$ws.Cells.Item(44, 3).Value = 0.96
$ws.Cells.Item(44, 4).Value = 0.96
$ws.Cells.Item(44, 5).Value = 0.96
$ws.Cells.Item(44, 6).Value = 0.96
$ws.Cells.Item(44, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(44, 8).Value = "v1.0.3"
$ws.Cells.Item(44, 9).Value = 2893
$ws.Cells.Item(44, 10).Value = 0.5808
$ws.Cells.Item(44, 11).Value = 0.96
$ws.Cells.Item(44, 12).Value = 0.2904
$ws.Cells.Item(44, 13).Value = 0.48
$ws.Cells.Item(44, 14).Value = 0.7704

$ws.Cells.Item(45, 1).Value = 0.7884615384615384
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = 0.9230769230769231
$ws.Cells.Item(45, 4).Value = 0.9230769230769231
$ws.Cells.Item(45, 5).Value = 0.9230769230769231
$ws.Cells.Item(45, 6).Value = 0.9230769230769231
$ws.Cells.Item(45, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(45, 8).Value = "v1.0.4"
$ws.Cells.Item(45, 9).Value = 3057
$ws.Cells.Item(45, 10).Value = 0.5648076923076923
$ws.Cells.Item(45, 11).Value = 0.9230769230769231
$ws.Cells.Item(45, 12).Value = 0.2824038461538462
$ws.Cells.Item(45, 13).Value = 0.4615384615384616
$ws.Cells.Item(45, 14).Value = 0.7439423076923077

$ws.Cells.Item(46, 1).Value = 0.7692307692307693
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(46, 3).Value = 0.9230769230769231
$ws.Cells.Item(46, 4).Value = 0.9230769230769231
$ws.Cells.Item(46, 5).Value = 0.9230769230769231
$ws.Cells.Item(46, 6).Value = 0.9230769230769231
$ws.Cells.Item(46, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(46, 8).Value = "v1.0.5"
$ws.Cells.Item(46, 9).Value = 3091
$ws.Cells.Item(46, 10).Value = 0.5584615384615386
$ws.Cells.Item(46, 11).Value = 0.9230769230769231
$ws.Cells.Item(46, 12).Value = 0.2792307692307693
$ws.Cells.Item(46, 13).Value = 0.4615384615384616
$ws.Cells.Item(46, 14).Value = 0.7407692307692308

$ws.Cells.Item(47, 1).Value = 0.7692307692307693
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(47, 3).Value = 0.9230769230769231
$ws.Cells.Item(47, 4).Value = 0.9230769230769231
$ws.Cells.Item(47, 5).Value = 0.9230769230769231
$ws.Cells.Item(47, 6).Value = 0.9230769230769231
$ws.Cells.Item(47, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(47, 8).Value = "v1.1.0"
$ws.Cells.Item(47, 9).Value = 3091
$ws.Cells.Item(47, 10).Value = 0.5584615384615386
$ws.Cells.Item(47, 11).Value = 0.9230769230769231
$ws.Cells.Item(47, 12).Value = 0.2792307692307693
$ws.Cells.Item(47, 13).Value = 0.4615384615384616
$ws.Cells.Item(47, 14).Value = 0.7407692307692308

$ws.Cells.Item(48, 1).Value = 0.7857142857142857
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(48, 3).Value = 0.9285714285714286
$ws.Cells.Item(48, 4).Value = 0.9285714285714286
$ws.Cells.Item(48, 5).Value = 0.9285714285714286
$ws.Cells.Item(48, 6).Value = 0.9285714285714286
$ws.Cells.Item(48, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(48, 8).Value = "v1.1.1"
$ws.Cells.Item(48, 9).Value = 3545
$ws.Cells.Item(48, 10).Value = 0.5657142857142857
$ws.Cells.Item(48, 11).Value = 0.9285714285714286
$ws.Cells.Item(48, 12).Value = 0.2828571428571429
$ws.Cells.Item(48, 13).Value = 0.4642857142857143
$ws.Cells.Item(48, 14).Value = 0.7471428571428571

$ws.Cells.Item(49, 1).Value = 0.7857142857142857
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(49, 3).Value = 0.9285714285714286
$ws.Cells.Item(49, 4).Value = 0.9285714285714286
$ws.Cells.Item(49, 5).Value = 0.9285714285714286
$ws.Cells.Item(49, 6).Value = 0.9285714285714286
$ws.Cells.Item(49, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(49, 8).Value = "v1.1.2"
$ws.Cells.Item(49, 9).Value = 3511
$ws.Cells.Item(49, 10).Value = 0.5657142857142857
$ws.Cells.Item(49, 11).Value = 0.9285714285714286
$ws.Cells.Item(49, 12).Value = 0.2828571428571429
$ws.Cells.Item(49, 13).Value = 0.4642857142857143
$ws.Cells.Item(49, 14).Value = 0.7471428571428571

$ws.Cells.Item(50, 1).Value = 0.7857142857142857
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(50, 3).Value = 0.9285714285714286
$ws.Cells.Item(50, 4).Value = 0.9285714285714286
$ws.Cells.Item(50, 5).Value = 0.9285714285714286
$ws.Cells.Item(50, 6).Value = 0.9285714285714286
$ws.Cells.Item(50, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(50, 8).Value = "v2.0.0"
$ws.Cells.Item(50, 9).Value = 3511
$ws.Cells.Item(50, 10).Value = 0.5657142857142857
$ws.Cells.Item(50, 11).Value = 0.9285714285714286
$ws.Cells.Item(50, 12).Value = 0.2828571428571429
$ws.Cells.Item(50, 13).Value = 0.4642857142857143
$ws.Cells.Item(50, 14).Value = 0.7471428571428571

$ws.Cells.Item(51, 1).Value = 0.7796610169491526
$ws.Cells.Item(51, 2).Value = 0.01694915254237288
$ws.Cells.Item(51, 3).Value = 0.9322033898305084
$ws.Cells.Item(51, 4).Value = 0.9322033898305084
$ws.Cells.Item(51, 5).Value = 0.9322033898305084
$ws.Cells.Item(51, 6).Value = 0.9322033898305084
$ws.Cells.Item(51, 7).Value = "Cartografia-social-api-front"
$ws.Cells.Item(51, 8).Value = "v2.0.1"
$ws.Cells.Item(51, 9).Value = 3695
$ws.Cells.Item(51, 10).Value = 0.5705084745762712
$ws.Cells.Item(51, 11).Value = 0.9322033898305084
$ws.Cells.Item(51, 12).Value = 0.2852542372881356
$ws.Cells.Item(51, 13).Value = 0.4661016949152542
$ws.Cells.Item(51, 14).Value = 0.7513559322033898
